$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.172.20"
$ws.Range("D3").Value = "2.996.65"
$ws.Range("E3").Value = "  +3.02%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'580.74"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "'162.78"
$ws.Range("E6").Value = "  +12.57%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.13%  "
$ws.Range("D9").Value = "2.993.43"
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("D10").Value = "'6.49"
$ws.Range("E10").Value = "  -6.06%  "
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  +3.03%  "
$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = "  +4.88%  "
$ws.Range("E13").Value = "  +5.60%  "
$ws.Range("D14").Value = "'34.51"
$ws.Range("E14").Value = "  +4.79%  "
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").Value = "66.204.56"
$ws.Range("E16").Value = "  +5.66%  "
$ws.Range("D17").Value = "3.495.69"
$ws.Range("D18").Value = "'6.90"
$ws.Range("D19").Value = "2.997.54"
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").Value = "'452.85"
$ws.Range("E20").Value = "  +5.28%  "
$ws.Range("D21").Value = "'13.82"
$ws.Range("E21").Value = "  +5.12%  "
$ws.Range("D22").Value = "'0.685"
$ws.Range("E22").Value = "  +3.67%  "
$ws.Range("D23").Value = "'7.32"
$ws.Range("E23").Value = "  +6.00%  "
$ws.Range("D24").Value = "'82.17"
$ws.Range("E24").Value = "  +4.23%  "
$ws.Range("E25").Value = "  +13.56%  "
$ws.Range("E26").Value = "  +2.17%  "
$ws.Range("D27").Value = "'10.37"
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").Value = "'8.12"
$ws.Range("E29").Value = "  +12.73%  "
$ws.Range("E30").Value = "  +18.25%  "
$ws.Range("E31").Value = "  +4.77%  "
$ws.Range("E32").Value = "  -6.42%  "
$ws.Range("D33").Value = "'27.20"
$ws.Range("E33").Value = "  +5.07%  "
$ws.Range("E34").Value = "  +3.05%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "'0.989"
$ws.Range("E36").Value = "  +3.55%  "
$ws.Range("E37").Value = "  +7.46%  "
$ws.Range("D38").Value = "'2.05"
$ws.Range("E38").Value = "  +7.30%  "
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").Value = "'2.96"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +15.36%  "
$ws.Range("E42").Value = "  +6.71%  "
$ws.Range("D43").Value = "'43.75"
$ws.Range("E43").Value = "  +6.35%  "
$ws.Range("E44").Value = "  +3.86%  "
$ws.Range("D45").Value = "'397.97"
$ws.Range("E45").Value = "  +11.15%  "
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("D47").Value = "2.786.95"
$ws.Range("E47").Value = "  +2.59%  "
$ws.Range("D48").Value = "'133.29"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D50").Value = "'23.79"
$ws.Range("E50").Value = "  +10.76%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  +3.28%  "
